$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8-17, columns B (name), C (from_bus), D (to_bus), E (in_service)
# Row 8 and 9 become the new "line7"/"line8" entries; rows that used to be
# extr1..extr8 (old rows 8-15) shift down to rows 10-17, and gain two more
# (extr7, extr8) that had previously been the last row but now occupy new
# rows 16 and 17.

$data = @(
    @(8,  "line7", 14, 11, $true),
    @(9,  "line8", 16,  9, $true),
    @(10, "extr1",  5, 12, $true),
    @(11, "extr2",  5,  9, $true),
    @(12, "extr3", 10, 11, $true),
    @(13, "extr4",  7,  8, $false),
    @(14, "extr5",  9, 11, $true),
    @(15, "extr6",  7, 11, $true),
    @(16, "extr7",  5,  7, $true),
    @(17, "extr8",  8,  5, $true)
)

# Reference cell carrying the "A-column" style (bold font, thin box border,
# center/top alignment) that every data row in column A uses.
$styleSrc = $ws.Cells.Item(2, 1)

foreach ($item in $data) {
    $r = $item[0]
    $name = $item[1]
    $fromBus = $item[2]
    $toBus = $item[3]
    $inService = $item[4]

    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value = $r - 2
    $aCell.Font.Bold = $styleSrc.Font.Bold
    $aCell.HorizontalAlignment = $styleSrc.HorizontalAlignment
    $aCell.VerticalAlignment = $styleSrc.VerticalAlignment
    $aCell.Borders.LineStyle = $styleSrc.Borders.LineStyle

    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $fromBus
    $ws.Cells.Item($r, 4).Value = $toBus
    $ws.Cells.Item($r, 5).Value = $inService
}
